$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("B2").Value = "Thu Jun 19 18:35:31 IST 2025"
$ws.Range("C2").Value = "Pass"

$ws = $wb.Worksheets.Item(2)
$ws.Range("B2").Value = "Thu Jun 19 18:40:14 IST 2025"
$ws.Range("C2").Value = "Pass"
$ws.Range("B3").Value = "Thu Jun 19 18:41:14 IST 2025"
$ws.Range("C3").Value = "Pass"
$ws.Range("B4").Value = "Thu Jun 19 18:42:11 IST 2025"
$ws.Range("C4").Value = "Pass"
$ws.Range("B5").Value = "Thu Jun 19 18:43:02 IST 2025"
$ws.Range("C5").Value = "Pass"

$ws = $wb.Worksheets.Item(3)
$ws.Range("B2").Value = "Thu Jun 19 18:36:27 IST 2025"
$ws.Range("C2").Value = "Pass"
$ws.Range("B3").Value = "Thu Jun 19 18:37:25 IST 2025"
$ws.Range("C3").Value = "Pass"
$ws.Range("B4").Value = "Thu Jun 19 18:38:19 IST 2025"
$ws.Range("C4").Value = "Pass"
$ws.Range("B5").Value = "Thu Jun 19 18:39:19 IST 2025"
$ws.Range("C5").Value = "Pass"

$ws = $wb.Worksheets.Item(4)
$ws.Range("B2").Value = "Thu Jun 19 18:24:13 IST 2025"
$ws.Range("C2").Value = "Pass"
$ws.Range("B3").Value = "Thu Jun 19 18:25:10 IST 2025"
$ws.Range("C3").Value = "Pass"
$ws.Range("B4").Value = "Thu Jun 19 18:26:04 IST 2025"
$ws.Range("C4").Value = "Pass"
$ws.Range("B5").Value = "Thu Jun 19 18:26:58 IST 2025"
$ws.Range("C5").Value = "Pass"

$ws = $wb.Worksheets.Item(5)
$ws.Range("B2").Value = "Thu Jun 19 18:31:53 IST 2025"
$ws.Range("C2").Value = "Pass"
$ws.Range("B3").Value = "Thu Jun 19 18:32:44 IST 2025"
$ws.Range("C3").Value = "Pass"
$ws.Range("B4").Value = "Thu Jun 19 18:33:38 IST 2025"
$ws.Range("C4").Value = "Pass"
$ws.Range("B5").Value = "Thu Jun 19 18:34:32 IST 2025"
$ws.Range("C5").Value = "Pass"

$ws = $wb.Worksheets.Item(6)
$ws.Range("B2").Value = "Thu Jun 19 18:27:57 IST 2025"
$ws.Range("C2").Value = "Pass"
$ws.Range("B3").Value = "Thu Jun 19 18:29:07 IST 2025"
$ws.Range("C3").Value = "Pass"
$ws.Range("B4").Value = "Thu Jun 19 18:30:09 IST 2025"
$ws.Range("C4").Value = "Pass"
$ws.Range("B5").Value = "Thu Jun 19 18:31:00 IST 2025"
$ws.Range("C5").Value = "Pass"

$ws = $wb.Worksheets.Item(7)
$ws.Range("B2").Value = "Thu Jun 19 18:21:42 IST 2025"
$ws.Range("C2").Value = "Pass"

$ws = $wb.Worksheets.Item(8)
$ws.Range("B2").Value = "Thu Jun 19 18:23:20 IST 2025"
$ws.Range("C2").Value = "Pass"

$ws = $wb.Worksheets.Item(9)
$ws.Range("B2").Value = "Thu Jun 19 18:22:32 IST 2025"
$ws.Range("C2").Value = "Pass"
